$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (within the used range), shifting existing data right to E:L
$ws.Range("D5:D102").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# Copy formatting (number format/style) from column E (which now holds the former column D data) into the new column D
$ws.Range("E5:E102").Copy($ws.Range("D5:D102"))

# Populate the new column D with the new period (fiscal year ending 2018-12-29) figures
$ws.Range("D7").Value = 43463
$ws.Range("D8").Value = 1718500
$ws.Range("D9").Value = 1065900
$ws.Range("D10").Value = 652500
$ws.Range("D12").Value = 87300
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 52200
$ws.Range("D17").Value = 1493400
$ws.Range("D18").Value = 225000
$ws.Range("D20").Value = 2500
$ws.Range("D21").Value = 330700
$ws.Range("D22").Value = 22600
$ws.Range("D23").Value = 204900
$ws.Range("D24").Value = 37100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 167800
$ws.Range("D27").Value = 167800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -3200
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2500
$ws.Range("D33").Value = 164600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 164600
$ws.Range("D38").Value = 43463
$ws.Range("D41").Value = 489700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 235200
$ws.Range("D44").Value = 258200
$ws.Range("D45").Value = 49300
$ws.Range("D46").Value = 1032500
$ws.Range("D47").Value = 25400
$ws.Range("D48").Value = 339900
$ws.Range("D49").Value = 1188200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 28300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2614300
$ws.Range("D57").Value = 126300
$ws.Range("D58").Value = 10000
$ws.Range("D59").Value = 159000
$ws.Range("D60").Value = 295300
$ws.Range("D61").Value = 684700
$ws.Range("D62").Value = 156000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1136100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 856500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1478200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("D81").Value = 164600
$ws.Range("D83").Value = 103200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 331800
$ws.Range("D91").Value = -74800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -382200
$ws.Range("D96").Value = -40000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 121900
$ws.Range("D101").Value = -11400
$ws.Range("D102").Value = 60100
